$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 204
$wsExpo.Range("F4").Value = 810
$wsExpo.Range("F6").Value = 23

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 204
$wsAll.Range("F5").Value = 810
$wsAll.Range("F7").Value = 23
